$wb = $excel.ActiveWorkbook

# The "完成量" (Q) column for sheet "1-1" was filled with the literal text
# "NaN" for rows 6-23 and 25-27, which made the dependent R/S percentage
# and MAX() formulas evaluate to #VALUE! errors. Replace those placeholder
# text values with the real numeric completion value (1) so the sheet
# recalculates cleanly, matching sheets "1-2"/"1-3"/"1-4" which already use
# numeric values in the same column.
$ws = $wb.Worksheets.Item("1-1")
$ws.Range("Q6:Q23").Value = 1
$ws.Range("Q25:Q27").Value = 1
